# Generate Report for handback
# Advance the handoff/handback timestamps for the first (16be212c...) file
# row on both the zh-cn and de-de report sheets, reflecting a fresh
# handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 07:20:45"
$wsZhCn.Range("G2").Value = "2016-01-08 07:21:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 07:20:55"
$wsDeDe.Range("G2").Value = "2016-01-08 07:21:50"
